$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0
$ws.Range("F2").Value = 25.46000000000054
$ws.Range("H2").Value = 0.07744482676076225
$ws.Range("I2").Value = 0.07744482676076225
$ws.Range("L2").Value = 6.073206112278585
$ws.Range("M2").Value = '[-0.5472695742418505, 12.693681798799021]'
$ws.Range("N2").Value = 0.07123838904066093
$ws.Range("O2").Value = 0.07123838904066093
$ws.Range("P2").Value = -1.434000250287233
$ws.Range("Q2").Value = '[-2.968632097085851, 0.10063159651138509]'
$ws.Range("R2").Value = 0.06630850168458191
$ws.Range("S2").Value = 0.06630850168458191
$ws.Range("T2").Value = 9.484362363770115
$ws.Range("U2").Value = '[5.650081247196265, 13.318643480343965]'
$ws.Range("V2").Value = 0.000009766253963539029
$ws.Range("W2").Value = 0.000009766253963539029
$ws.Range("X2").Value = 5.810690690690816
$ws.Range("Y2").Value = -0.4077677677677753
$ws.Range("Z2").Value = 12.02914914914941
$ws.Range("F3").Value = 25.46000000000054
$ws.Range("H3").Value = 0.07599240397497653
$ws.Range("I3").Value = 0.07599240397497653
$ws.Range("L3").Value = 6.358506434937368
$ws.Range("M3").Value = '[-1.1135990494853782, 13.830611919360114]'
$ws.Range("N3").Value = 0.09342240760797393
$ws.Range("O3").Value = 0.09342240760797393
$ws.Range("P3").Value = -1.081789662497386
$ws.Range("Q3").Value = '[-3.0566847440333125, 0.8931054190385397]'
$ws.Range("R3").Value = 0.2757764751789915
$ws.Range("S3").Value = 0.2757764751789915
$ws.Range("T3").Value = 10.57418878981332
$ws.Range("U3").Value = '[6.584643374187706, 14.563734205438934]'
$ws.Range("V3").Value = 0.000002959334910679345
$ws.Range("W3").Value = 0.000002959334910679345
$ws.Range("X3").Value = 4.383503503503597
$ws.Range("Y3").Value = -3.618938938939015
$ws.Range("Z3").Value = 12.38594594594621
$ws.Range("F4").Value = 25.46000000000054
$ws.Range("H4").Value = 0.07995015734492372
$ws.Range("I4").Value = 0.07995015734492372
$ws.Range("L4").Value = 6.142802739779688
$ws.Range("M4").Value = '[-0.3467085606095015, 12.632314040168877]'
$ws.Range("N4").Value = 0.06298340929045554
$ws.Range("O4").Value = 0.06298340929045554
$ws.Range("P4").Value = -1.58494764505431
$ws.Range("Q4").Value = '[-2.993789996213697, -0.17610529389492324]'
$ws.Range("R4").Value = 0.02831859753591659
$ws.Range("S4").Value = 0.02831859753591659
$ws.Range("T4").Value = 10.2195985060718
$ws.Range("U4").Value = '[6.324684770720216, 14.114512241423393]'
$ws.Range("V4").Value = 0.000003545816271932622
$ws.Range("W4").Value = 0.000003545816271932622
$ws.Range("X4").Value = 6.422342342342478
$ws.Range("Y4").Value = 0.7135935935936084
$ws.Range("Z4").Value = 12.13109109109135
$ws.Range("F5").Value = 25.46000000000054
$ws.Range("H5").Value = 0.05839504112816929
$ws.Range("I5").Value = 0.05839504112816929
$ws.Range("L5").Value = 7.462490555871049
$ws.Range("M5").Value = '[-0.649565942822047, 15.574547054564146]'
$ws.Range("N5").Value = 0.07047211424636068
$ws.Range("O5").Value = 0.07047211424636068
$ws.Range("P5").Value = -2.012631930227695
$ws.Range("Q5").Value = '[-3.761105919613004, -0.2641579408423853]'
$ws.Range("R5").Value = 0.02502949273132815
$ws.Range("S5").Value = 0.02502949273132815
$ws.Range("T5").Value = 11.52592157725771
$ws.Range("U5").Value = '[7.174739762245093, 15.877103392270332]'
$ws.Range("V5").Value = 0.000002990797226987141
$ws.Range("W5").Value = 0.000002990797226987141
$ws.Range("X5").Value = 8.155355355355528
$ws.Range("Y5").Value = 1.070390390390417
$ws.Range("Z5").Value = 15.24032032032064
$ws.Range("F6").Value = 25.46000000000054
$ws.Range("H6").Value = 0.1022658604930162
$ws.Range("I6").Value = 0.1022658604930162
$ws.Range("L6").Value = 6.354634521749102
$ws.Range("M6").Value = '[-1.1854100571813078, 13.894679100679511]'
$ws.Range("N6").Value = 0.09651727250618092
$ws.Range("O6").Value = 0.09651727250618092
$ws.Range("P6").Value = -2.452895164965004
$ws.Range("Q6").Value = '[-5.559895707254008, 0.6541053773240009]'
$ws.Range("R6").Value = 0.1188182318997735
$ws.Range("S6").Value = 0.1188182318997735
$ws.Range("T6").Value = 10.59033264056075
$ws.Range("U6").Value = '[6.537847102621363, 14.642818178500136]'
$ws.Range("V6").Value = 0.000003808359848811449
$ws.Range("W6").Value = 0.000003808359848811449
$ws.Range("X6").Value = 9.939339339339551
$ws.Range("Y6").Value = -2.650490490490546
$ws.Range("Z6").Value = 22.52916916916965
$ws.Range("F7").Value = 25.46000000000054
$ws.Range("H7").Value = 0.1467725149145001
$ws.Range("I7").Value = 0.1467725149145001
$ws.Range("L7").Value = 5.05588746357023
$ws.Range("M7").Value = '[-1.5309532078364363, 11.642728134976897]'
$ws.Range("N7").Value = 0.1291154472078895
$ws.Range("O7").Value = 0.1291154472078895
$ws.Range("P7").Value = -2.742211004935235
$ws.Range("Q7").Value = '[-5.15736932120847, -0.32705268866200043]'
$ws.Range("R7").Value = 0.02696124826939617
$ws.Range("S7").Value = 0.02696124826939617
$ws.Range("T7").Value = 9.593399832124881
$ws.Range("U7").Value = '[6.100202718012323, 13.086596946237439]'
$ws.Range("V7").Value = 0.000001540700505398007
$ws.Range("W7").Value = 0.000001540700505398007
$ws.Range("X7").Value = 11.11167167167191
$ws.Range("Y7").Value = 1.325245245245275
$ws.Range("Z7").Value = 20.89809809809854
$ws.Range("F8").Value = 25.46000000000054
$ws.Range("H8").Value = 0.09362642948592337
$ws.Range("I8").Value = 0.09362642948592337
$ws.Range("L8").Value = 6.479935514015926
$ws.Range("M8").Value = '[-0.8199047874391052, 13.779775815470957]'
$ws.Range("N8").Value = 0.08052980491501227
$ws.Range("O8").Value = 0.08052980491501227
$ws.Range("Q8").Value = '[-4.528421843012315, -1.2327370572644627]'
$ws.Range("R8").Value = 0.0009982735354077121
$ws.Range("S8").Value = 0.0009982735354077121
$ws.Range("T8").Value = 10.44105934797114
$ws.Range("U8").Value = '[6.438388132482297, 14.443730563459985]'
$ws.Range("V8").Value = 0.000003933440194003879
$ws.Range("W8").Value = 0.000003933440194003879
$ws.Range("X8").Value = 11.6723523523526
$ws.Range("Y8").Value = 4.995155155155257
$ws.Range("Z8").Value = 18.34954954954994
$ws.Range("F9").Value = 23.96000000000031
$ws.Range("H9").Value = 0.08016064793713551
$ws.Range("I9").Value = 0.08016064793713551
$ws.Range("L9").Value = 6.184209720110208
$ws.Range("M9").Value = '[-0.777579920502733, 13.145999360723149]'
$ws.Range("N9").Value = 0.08032374259788266
$ws.Range("O9").Value = 0.08032374259788266
$ws.Range("P9").Value = 2.584974135386196
$ws.Range("Q9").Value = '[0.19497371824080822, 4.9749745525315845]'
$ws.Range("R9").Value = 0.03465405565016577
$ws.Range("S9").Value = 0.03465405565016577
$ws.Range("T9").Value = 9.737588378453177
$ws.Range("U9").Value = '[5.927242762285392, 13.547933994620962]'
$ws.Range("V9").Value = 0.000005626742517694794
$ws.Range("W9").Value = 0.000005626742517694794
$ws.Range("X9").Value = 14.10258258258276
$ws.Range("Y9").Value = 4.988668668668732
$ws.Range("Z9").Value = 23.2164964964968
$ws.Range("F10").Value = 23.96000000000031
$ws.Range("H10").Value = 0.1217704737103998
$ws.Range("I10").Value = 0.1217704737103998
$ws.Range("L10").Value = 6.019530825406544
$ws.Range("M10").Value = '[-1.4217782981584381, 13.460839948971525]'
$ws.Range("N10").Value = 0.1102357380339432
$ws.Range("O10").Value = 0.1102357380339432
$ws.Range("P10").Value = 2.912026824048196
$ws.Range("Q10").Value = '[1.1509738850989621, 4.67307976299743]'
$ws.Range("R10").Value = 0.001738160074938877
$ws.Range("S10").Value = 0.001738160074938877
$ws.Range("T10").Value = 10.30130048424756
$ws.Range("U10").Value = '[6.22115744383561, 14.3814435246595]'
$ws.Range("V10").Value = 0.000006925884731723286
$ws.Range("W10").Value = 0.000006925884731723286
$ws.Range("X10").Value = 12.85541541541558
$ws.Range("Y10").Value = 6.139899899899975
$ws.Range("Z10").Value = 19.57093093093118
$ws.Range("F11").Value = 23.96000000000031
$ws.Range("H11").Value = 0.1264210264487247
$ws.Range("I11").Value = 0.1264210264487247
$ws.Range("L11").Value = 5.568826615557604
$ws.Range("M11").Value = '[-1.425266596868937, 12.562919827984146]'
$ws.Range("N11").Value = 0.1157842652413437
$ws.Range("O11").Value = 0.1157842652413437
$ws.Range("P11").Value = 2.949763672739966
$ws.Range("Q11").Value = '[1.2264475824825016, 4.67307976299743]'
$ws.Range("R11").Value = 0.00123812719928762
$ws.Range("S11").Value = 0.00123812719928762
$ws.Range("T11").Value = 10.61315395008369
$ws.Range("U11").Value = '[6.803102993987062, 14.423204906180315]'
$ws.Range("V11").Value = 0.000001178109537702809
$ws.Range("W11").Value = 0.000001178109537702809
$ws.Range("X11").Value = 12.71151151151168
$ws.Range("Y11").Value = 6.139899899899979
$ws.Range("Z11").Value = 19.28312312312337
$ws.Range("F12").Value = 23.96000000000031
$ws.Range("H12").Value = 0.09134450782756509
$ws.Range("I12").Value = 0.09134450782756509
$ws.Range("L12").Value = 6.576092087862708
$ws.Range("M12").Value = '[-1.3153579633704275, 14.467542139095844]'
$ws.Range("N12").Value = 0.1002035557394314
$ws.Range("O12").Value = 0.1002035557394314
$ws.Range("P12").Value = 2.673026782333658
$ws.Range("Q12").Value = '[0.6352369529781159, 4.710816611689199]'
$ws.Range("R12").Value = 0.01129718079555841
$ws.Range("S12").Value = 0.01129718079555841
$ws.Range("T12").Value = 10.84281234259418
$ws.Range("U12").Value = '[6.674235705117968, 15.01138898007039]'
$ws.Range("V12").Value = 0.000004136690901157536
$ws.Range("W12").Value = 0.000004136690901157536
$ws.Range("X12").Value = 13.76680680680698
$ws.Range("Y12").Value = 5.995995995996074
$ws.Range("Z12").Value = 21.5376176176179
$ws.Range("F13").Value = 23.96000000000031
$ws.Range("H13").Value = 0.1223517643631069
$ws.Range("I13").Value = 0.1223517643631069
$ws.Range("L13").Value = 5.959739813379944
$ws.Range("M13").Value = '[-1.7794514773019845, 13.698931104061872]'
$ws.Range("N13").Value = 0.1279054862777846
$ws.Range("O13").Value = 0.1279054862777846
$ws.Range("P13").Value = 2.547237286694427
$ws.Range("Q13").Value = '[-0.5786316799404627, 5.673106253329316]'
$ws.Range("R13").Value = 0.1077129706982298
$ws.Range("S13").Value = 0.1077129706982298
$ws.Range("T13").Value = 9.23701976195867
$ws.Range("U13").Value = '[5.194093684697659, 13.279945839219682]'
$ws.Range("V13").Value = 0.00003420327495740061
$ws.Range("W13").Value = 0.00003420327495740061
$ws.Range("X13").Value = 14.24648648648667
$ws.Range("Y13").Value = 2.326446446446477
$ws.Range("Z13").Value = 26.16652652652686
$ws.Range("F14").Value = 23.96000000000031
$ws.Range("H14").Value = 0.09316998075902727
$ws.Range("I14").Value = 0.09316998075902727
$ws.Range("L14").Value = 5.764232383277372
$ws.Range("M14").Value = '[-0.7275282820354185, 12.255993048590163]'
$ws.Range("N14").Value = 0.08044784203367961
$ws.Range("O14").Value = 0.08044784203367961
$ws.Range("P14").Value = 3.012658420559581
$ws.Range("Q14").Value = '[1.452868674633116, 4.572448166486046]'
$ws.Range("R14").Value = 0.0003277309831815067
$ws.Range("S14").Value = 0.0003277309831815067
$ws.Range("T14").Value = 10.00943168956243
$ws.Range("U14").Value = '[6.338554596514679, 13.68030878261018]'
$ws.Range("V14").Value = 0.000001761171098513259
$ws.Range("W14").Value = 0.000001761171098513259
$ws.Range("X14").Value = 12.47167167167183
$ws.Range("Y14").Value = 6.523643643643723
$ws.Range("Z14").Value = 18.41969969969994
$ws.Range("F15").Value = 23.96000000000031
$ws.Range("H15").Value = 0.06308996557909863
$ws.Range("I15").Value = 0.06308996557909863
$ws.Range("L15").Value = 6.845807911508778
$ws.Range("M15").Value = '[0.2535691329220491, 13.438046690095508]'
$ws.Range("N15").Value = 0.0421511427206136
$ws.Range("O15").Value = 0.0421511427206136
$ws.Range("P15").Value = 3.138447916198812
$ws.Range("Q15").Value = '[1.8050792624229643, 4.47181656997466]'
$ws.Range("R15").Value = 0.00002169444178856672
$ws.Range("S15").Value = 0.00002169444178856672
$ws.Range("T15").Value = 10.84474082708262
$ws.Range("U15").Value = '[6.8220996824707765, 14.867381971694458]'
$ws.Range("V15").Value = 0.000002172409638179218
$ws.Range("W15").Value = 0.000002172409638179218
$ws.Range("X15").Value = 11.99199199199215
$ws.Range("Y15").Value = 6.907387387387478
$ws.Range("Z15").Value = 17.07659659659681
$ws.Range("F16").Value = 23.96000000000031
$ws.Range("H16").Value = 0.1919801547151072
$ws.Range("I16").Value = 0.1919801547151072
$ws.Range("L16").Value = 5.320733970098344
$ws.Range("M16").Value = '[-1.8926726997940602, 12.534140639990747]'
$ws.Range("N16").Value = 0.1443476159518164
$ws.Range("O16").Value = 0.1443476159518164
$ws.Range("P16").Value = 2.798816277972889
$ws.Range("Q16").Value = '[-0.32705268866200043, 5.924685244607778]'
$ws.Range("R16").Value = 0.07802747816666122
$ws.Range("S16").Value = 0.07802747816666122
$ws.Range("T16").Value = 10.86231334705414
$ws.Range("U16").Value = '[6.7680804908881775, 14.956546203220103]'
$ws.Range("V16").Value = 0.000002907490642023447
$ws.Range("W16").Value = 0.000002907490642023447
$ws.Range("X16").Value = 13.2871271271273
$ws.Range("Y16").Value = 1.367087087087105
$ws.Range("Z16").Value = 25.20716716716749
$ws.Range("F17").Value = 23.96000000000031
$ws.Range("H17").Value = 0.1960016336536936
$ws.Range("I17").Value = 0.1960016336536936
$ws.Range("L17").Value = 4.773930804930431
$ws.Range("M17").Value = '[-1.758429218348546, 11.306290828209407]'
$ws.Range("N17").Value = 0.1479988761163087
$ws.Range("O17").Value = 0.1479988761163087
$ws.Range("P17").Value = 2.937184723176043
$ws.Range("Q17").Value = '[0.03144737390980801, 5.842922072442278]'
$ws.Range("R17").Value = 0.04767273423852725
$ws.Range("S17").Value = 0.04767273423852725
$ws.Range("T17").Value = 10.86185598283731
$ws.Range("U17").Value = '[7.161499148117416, 14.562212817557201]'
$ws.Range("V17").Value = 0.0000004215925266670695
$ws.Range("W17").Value = 0.0000004215925266670695
$ws.Range("X17").Value = 12.75947947947964
$ws.Range("Y17").Value = 1.678878878878901
$ws.Range("Z17").Value = 23.84008008008038
